$wb = $excel.ActiveWorkbook

# The localization status changed from "Ready for handoff" to "In Translation"
# for the zh-cn / de-de columns, and the corresponding "Status" columns got
# narrower (their column width shrank from ~17.22 to ~13.41 character-units).
$newStatus = "In Translation"

# ColumnWidth (character units) gets quantized by the host when stored back
# to OOXML width. 12.5 is the closest input that lands on the nearest
# achievable stored width to the target 13.4101845877511.
$newColWidth = 12.5

# --- "Overview" sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- "zh-cn" detail sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- "de-de" detail sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
